$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the note in E61 ("Went up to here in updates")
$ws.Range("E61").ClearContents()

# Update detection points for the second half of the ground truth set
$updates = @{
    63  = @{ B = 29; C = 56 }
    64  = @{ B = 30; C = 49 }
    65  = @{ C = 50 }
    69  = @{ C = 32 }
    71  = @{ B = 29 }
    73  = @{ C = 51 }
    74  = @{ C = 47 }
    75  = @{ B = 32; C = 54 }
    76  = @{ C = 55 }
    77  = @{ C = 46 }
    78  = @{ C = 32 }
    80  = @{ C = 55 }
    82  = @{ C = 48 }
    83  = @{ B = 36; C = 32 }
    84  = @{ B = 33; C = 55 }
    87  = @{ B = 39; C = 47 }
    88  = @{ C = 32 }
    89  = @{ C = 48 }
    90  = @{ C = 53 }
    93  = @{ C = 38 }
    95  = @{ B = 41; C = 60 }
    97  = @{ C = 63 }
    98  = @{ B = 37; C = 52 }
    99  = @{ B = 23; C = 64 }
    101 = @{ B = 42; C = 52 }
    102 = @{ C = 48 }
    103 = @{ B = 28; C = 50 }
    105 = @{ B = 32; C = 27 }
    107 = @{ B = 27; C = 36 }
    108 = @{ B = 34; C = 56 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# Update the saved view state to match where the author left off
$ws.Application.ActiveWindow.ScrollRow = 87
$ws.Range("E96").Select()
